$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.398.90"
$ws.Range("E2").Value = "  -1.03%  "

# Row 3
$ws.Range("D3").Value = "3.473.20"
$ws.Range("E3").Value = "  -1.84%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.52"
$ws.Range("E5").Value = "  -1.91%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.00"
$ws.Range("E6").Value = "  -2.61%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.608"
$ws.Range("E8").Value = "  +2.18%  "

# Row 9
$ws.Range("D9").Value = "3.470.79"
$ws.Range("E9").Value = "  -1.81%  "

# Row 10
$ws.Range("E10").Value = "  -1.19%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.00"
$ws.Range("E11").Value = "  -2.74%  "

# Row 12
$ws.Range("E12").Value = "  -4.18%  "

# Row 13
$ws.Range("D13").Value = "4.078.58"
$ws.Range("E13").Value = "  -1.86%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.07"
$ws.Range("E14").Value = "  -1.08%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.133"
$ws.Range("E15").Value = "  -2.60%  "

# Row 16
$ws.Range("D16").Value = "67.404.19"
$ws.Range("E16").Value = "  -0.95%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000176"
$ws.Range("E17").Value = "  -2.70%  "

# Row 18
$ws.Range("D18").Value = "3.470.99"
$ws.Range("E18").Value = "  -2.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.12"
$ws.Range("E19").Value = "  -5.07%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.09"
$ws.Range("E20").Value = "  -4.57%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "388.26"
$ws.Range("E21").Value = "  -3.39%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.88"
$ws.Range("E22").Value = "  -2.46%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.79"
$ws.Range("E23").Value = "  +1.80%  "

# Row 24
$ws.Range("E24").Value = "  -0.09%  "

# Row 25
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.535"
$ws.Range("E25").Value = "  -2.47%  "

# Row 26
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.71"
$ws.Range("E26").Value = "  -3.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000120"
$ws.Range("E27").Value = "  -3.45%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.20"
$ws.Range("E28").Value = "  -5.61%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.176"
$ws.Range("E29").Value = "  -1.78%  "

# Row 30
$ws.Range("E30").Value = "  +0.32%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.06"
$ws.Range("E31").Value = "  -5.19%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.39"
$ws.Range("E32").Value = "  -5.79%  "

# Row 33
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.04"
$ws.Range("E33").Value = "  -2.38%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.05"
$ws.Range("E34").Value = "  -0.13%  "

# Row 35
$ws.Range("E35").Value = "  -0.12%  "

# Row 36
$ws.Range("E36").Value = "  -5.48%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.56"
$ws.Range("E37").Value = "  -5.78%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.86"
$ws.Range("E38").Value = "  -1.76%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.883"
$ws.Range("E39").Value = "  -0.27%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.73"
$ws.Range("E40").Value = "  -2.34%  "

# Row 41
$ws.Range("E41").Value = "  -5.80%  "

# Row 42
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.81"
$ws.Range("E42").Value = "  -1.39%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.62"
$ws.Range("E43").Value = "  -4.54%  "

# Row 44
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.56"
$ws.Range("E44").Value = "  -4.58%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0708"
$ws.Range("E45").Value = "  -5.03%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.03"
$ws.Range("E46").Value = "  -3.81%  "

# Row 47
$ws.Range("D47").Value = "2.716.04"
$ws.Range("E47").Value = "  -7.08%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.25"
$ws.Range("E48").Value = "  -3.06%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0298"
$ws.Range("E49").Value = "  -3.44%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "331.81"
$ws.Range("E50").Value = "  -6.60%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.04"
$ws.Range("E51").Value = "  -4.82%  "
